$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 5.7
$ws.Range("R2").Value = 1.42
$ws.Range("U2").Value = 2.04
$ws.Range("V2").Value = 1.21
$ws.Range("W2").Value = 2.4
$ws.Range("Y2").Value = 21
$ws.Range("AI2").Value = 75

# Row 4 updates
$ws.Range("J4").Value = 4.6
$ws.Range("AI4").Value = 48
